# Add meeting participation data down through row 49 (48 members total)
# and refresh the "Participation Data" sheet's member IDs / participation
# types to match the newly expanded roster.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Participation Data")

# Sequential Member IDs for rows 2-49 (member #1 .. #48)
$memberIds = 2..49

# Participation Type values (shared-string text) for rows 2-49
$participationTypes = @(
    "Table Topic","Table Topic","Table Topic","Table Topic","Table Topic",
    "Table Topic","Table Topic","Table Topic","Table Topic","Table Topic",
    "Table Topic","Prepared Speech","Grammarian","Timer","AH Counter",
    "Vote Counter","General Evaluator","Table Topics Master","Toastmaster of the Evening",
    "Prepared Speech","Prepared Speech","Prepared Speech","Evaluation","Evaluation",
    "Evaluation","Evaluation","Table Topic","Prepared Speech","Evaluation",
    "General Evaluator","Table Topics Master","Toastmaster of the Evening","Grammarian","Timer",
    "Vote Counter","AH Counter","Table Topic","Prepared Speech","Evaluation",
    "General Evaluator","Table Topics Master","Toastmaster of the Evening","Grammarian","Timer",
    "Vote Counter","AH Counter","Table Topic","Prepared Speech"
)

$startRow = 2
for ($i = 0; $i -lt $memberIds.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $memberIds[$i]
    $ws.Cells.Item($row, 3).Value = $participationTypes[$i]
}

# Match the workbook's updated used range / selection
[void]$ws.Range("C2:C11").Select()

Write-Output "Participation data expanded to row 49"
